$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove all existing hyperlinks (engine-level, affects whole sheet) ---
$ws.Range("A1").Hyperlinks.Delete()

# --- Header row ---
$ws.Range("A1").Value = "Url"
$ws.Range("B1").Value = "HTTP port"
$ws.Range("C1").Value = "Coap Port"
$ws.Range("D1").Value = "Coap secure Port"
$ws.Range("E1").Value = "Short ID"

# --- Row 2: admin (new) ---
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = 8079
$ws.Range("C2").Value = 5681
$ws.Range("D2").Value = 5682
$ws.Range("E2").Value = 1

# --- Row 3: demo ---
$ws.Range("A3").Value = "demo"
$ws.Range("B3").Value = 8080
$ws.Range("C3").Value = 5683
$ws.Range("D3").Value = 5684

# --- Row 4: mehariclub ---
$ws.Range("A4").Value = "mehariclub"
$ws.Range("B4").Value = 8081
$ws.Range("C4").Value = 5685
$ws.Range("D4").Value = 5686
$ws.Range("E4").Value = 3

# --- Row 5: v-mti (new) ---
$ws.Range("A5").Value = "v-mti"
$ws.Range("B5").Value = 8082
$ws.Range("C5").Value = 5687
$ws.Range("D5").Value = 5688
$ws.Range("E5").Value = 2

# --- Row 6: tyva (new) ---
$ws.Range("A6").Value = "tyva"
$ws.Range("B6").Value = 8083
$ws.Range("C6").Value = 5689
$ws.Range("D6").Value = 5690
$ws.Range("E6").Value = 4

# --- Hyperlinks: add in the original relationship order so rIds line up ---
$ws.Hyperlinks.Add($ws.Range("A4"), "https://mehariclub.iot.sheeld.co/")
$ws.Range("A4").Style = "Lien hypertexte"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://demo.iot.sheeld.co/")
$ws.Range("A3").Style = "Lien hypertexte"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://admin.iot.sheeld.co/")
$ws.Range("A2").Style = "Lien hypertexte"

$ws.Hyperlinks.Add($ws.Range("A5"), "https://v-mti.iot.sheeld.co/")
$ws.Range("A5").Style = "Lien hypertexte"

$ws.Hyperlinks.Add($ws.Range("A6"), "https://tyva.iot.sheeld.co/")
$ws.Range("A6").Style = "Lien hypertexte"

# --- Selection matches target sheet view ---
$ws.Range("F9").Select() | Out-Null
